$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header labels (feature names + Intercept)
$ws.Cells.Item(1,1).Value = "f1"
$ws.Cells.Item(1,2).Value = "f2"
$ws.Cells.Item(1,3).Value = "f3"
$ws.Cells.Item(1,4).Value = "f4"
$ws.Cells.Item(1,5).Value = "f5"
$ws.Cells.Item(1,6).Value = "f6"
$ws.Cells.Item(1,7).Value = "f7"
$ws.Cells.Item(1,8).Value = "f8"
$ws.Cells.Item(1,9).Value = "f9"
$ws.Cells.Item(1,10).Value = "f10"
$ws.Cells.Item(1,11).Value = "f11"
$ws.Cells.Item(1,12).Value = "f12"
$ws.Cells.Item(1,13).Value = "f13"
$ws.Cells.Item(1,14).Value = "f14"
$ws.Cells.Item(1,15).Value = "Intercept"

# Row 2: data (shifted up from old row 1, first column dropped, Intercept=1 appended)
$ws.Cells.Item(2,1).Value = 0.10563894528931221
$ws.Cells.Item(2,2).Value = 0.10046774502926963
$ws.Cells.Item(2,3).Value = 0.03801752165048895
$ws.Cells.Item(2,4).Value = 0.020799912435065453
$ws.Cells.Item(2,5).Value = 0.039949013933146496
$ws.Cells.Item(2,6).Value = 0.05058647573804735
$ws.Cells.Item(2,7).Value = 0.028155693205916356
$ws.Cells.Item(2,8).Value = 0.003931769553945087
$ws.Cells.Item(2,9).Value = 0.015275764776695604
$ws.Cells.Item(2,10).Value = 0.0077297707465816354
$ws.Cells.Item(2,11).Value = 0.01829950670864804
$ws.Cells.Item(2,12).Value = 0.014162830989866797
$ws.Cells.Item(2,13).Value = 0.008314546777205112
$ws.Cells.Item(2,14).Value = 0.003924116920803731
$ws.Cells.Item(2,15).Value = 1

# Row 3: data (shifted up from old row 2, first column dropped, Intercept=1 appended)
$ws.Cells.Item(3,1).Value = 0.054966196053412464
$ws.Cells.Item(3,2).Value = 0.04456314078420359
$ws.Cells.Item(3,3).Value = 0.05614385178041718
$ws.Cells.Item(3,4).Value = 0.023484171583938087
$ws.Cells.Item(3,5).Value = 0.012127697361843822
$ws.Cells.Item(3,6).Value = 0.01414531971871711
$ws.Cells.Item(3,7).Value = 0.024317825849495007
$ws.Cells.Item(3,8).Value = 0.030546846457921686
$ws.Cells.Item(3,9).Value = 0.02372971493587985
$ws.Cells.Item(3,10).Value = 0.03237993904930274
$ws.Cells.Item(3,11).Value = 0.023976529530875652
$ws.Cells.Item(3,12).Value = 0.02372106496484171
$ws.Cells.Item(3,13).Value = 0.032799314988443774
$ws.Cells.Item(3,14).Value = 0.041740925619264235
$ws.Cells.Item(3,15).Value = 1

# Row 4: data (shifted up from old row 3, first column dropped, Intercept=1 appended)
$ws.Cells.Item(4,1).Value = 0.21044496538405794
$ws.Cells.Item(4,2).Value = 0.020256031214466817
$ws.Cells.Item(4,3).Value = 0.04306763979832916
$ws.Cells.Item(4,4).Value = 0.03620466927126364
$ws.Cells.Item(4,5).Value = 0.0513276004066082
$ws.Cells.Item(4,6).Value = 0.014147930756331365
$ws.Cells.Item(4,7).Value = 0.04847761704586361
$ws.Cells.Item(4,8).Value = 0.04122984020138278
$ws.Cells.Item(4,9).Value = 0.0198755785293201
$ws.Cells.Item(4,10).Value = 0.02472948281297303
$ws.Cells.Item(4,11).Value = 0.019170943265950494
$ws.Cells.Item(4,12).Value = 0.01688534344992194
$ws.Cells.Item(4,13).Value = 0.017668996577541573
$ws.Cells.Item(4,14).Value = 0.013460499775156958
$ws.Cells.Item(4,15).Value = 1

# Row 5: data (shifted up from old row 4, first column dropped, Intercept=1 appended)
$ws.Cells.Item(5,1).Value = 0.032048952257576334
$ws.Cells.Item(5,2).Value = 0.140927773487289
$ws.Cells.Item(5,3).Value = 0.06932186929181296
$ws.Cells.Item(5,4).Value = 0.0799163061127188
$ws.Cells.Item(5,5).Value = 0.046430101614712393
$ws.Cells.Item(5,6).Value = 0.06876901899526384
$ws.Cells.Item(5,7).Value = 0.03381745298902126
$ws.Cells.Item(5,8).Value = 0.055505273539062476
$ws.Cells.Item(5,9).Value = 0.05887096677779242
$ws.Cells.Item(5,10).Value = 0.04483340883336165
$ws.Cells.Item(5,11).Value = 0.0424446822922641
$ws.Cells.Item(5,12).Value = 0.046985400130180814
$ws.Cells.Item(5,13).Value = 0.03647005230650012
$ws.Cells.Item(5,14).Value = 0.025804356293747235
$ws.Cells.Item(5,15).Value = 1

# Update selection to match the new used range
[void]$ws.Range("A1:O5").Select()
